$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.72356903553009
$ws.Range("B1").Value = 2.265031576156616
$ws.Range("C1").Value = 4.6170654296875
$ws.Range("D1").Value = 4.165021419525146
$ws.Range("E1").Value = 1.59408438205719
